# ------------------------------------------------------------------
# Helper functions
# ------------------------------------------------------------------
function Set-RowValues {
    param($ws, $row, $startCol, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $v = $values[$i]
        if ($v -ne $null) {
            $col = $startCol + $i
            $ws.Cells.Item($row, $col).Value = $v
        }
    }
}

function Set-TextValue {
    param($ws, $row, $col, $value)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "Đơn sale chính" - add two new order rows and recompute
# the "Tổng" (total) row, which moves from row 4 down to row 6.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# The old row 4 ("Tổng") is being replaced by fresh detail rows, and the
# recomputed "Tổng" row moves down to row 6 - clear any leftover cell
# content first so stale values from the old layout cannot survive.
$ws1.Range("A4:N20").ClearContents()

Set-RowValues $ws1 4 1 @("HD-LUXURY", 627, $null, "SÓC TRĂNG", "tạ duy hoàng ", "Cá nhân", "Cắt mí", 6000000, $null, $null, 6000000, 6000000, 0.13, 780000)
Set-TextValue $ws1 4 3 "08-04-2024"

Set-RowValues $ws1 5 1 @("HD-LUXURY", 628, $null, "SÓC TRĂNG", "nguyễn thị lệ trang", "Cá nhân", "Cắt mí", 4000000, $null, $null, 4000000, 4000000, 0.1, 400000)
Set-TextValue $ws1 5 3 "08-04-2024"

Set-RowValues $ws1 6 1 @("Tổng", 4, $null, $null, $null, $null, $null, 24000000, $null, 0, 24000000, 22000000, 0, 2560000)

# ------------------------------------------------------------------
# Sheet order / names: insert a new "Đơn thu nợ" sheet between
# "Đơn sale chính" and "Lương", and rebuild "Lương" so the sheetId
# sequence (1,2,3) matches the recorded edit exactly.
# ------------------------------------------------------------------
$wsLuongOld = $wb.Worksheets.Item(2)
$wsLuongOld.Delete()

$wsDebt = $wb.Worksheets.Add($null, $ws1)
$wsDebt.Name = "Đơn thu nợ"

$wsLuong = $wb.Worksheets.Add($null, $wsDebt)
$wsLuong.Name = "Lương"

# ------------------------------------------------------------------
# Sheet 2: "Đơn thu nợ" (brand new sheet)
# ------------------------------------------------------------------
Set-RowValues $wsDebt 1 1 @("Tiền tố", "Mã đơn thu nợ", "Lượng thu", "Ngày thu", "Cơ sở", "Đơn nợ", "Tên dịch vụ", "Khách hàng", "Nguồn khách", "Sale chính", "Đơn giá gốc", "Sale phụ", "Upsale", "Đơn giá", "Đã thanh toán", "Bác sĩ 1", "Bác sĩ 2", "Tỉ lệ chiết khấu sale chính", "Chiết khấu sale chính", "Tỉ lệ chiết khấu sale phụ", "Chiết khấu sale phụ", "Tỉ lệ chiết khấu bác sĩ 1", "Chiết khấu bác sĩ 1", "Tỉ lệ chiết khấu bác sĩ 2", "Chiết khấu bác sĩ 2")

Set-RowValues $wsDebt 2 1 @("TN", 177, 8000000, $null, "SÓC TRĂNG", "HD-LUXURY-611", "Nâng mũi", "ngọc hân", "CTV", "Thạch Hoàng Nhân", 35000000, "Lê Đình Hậu", 8000000, 43000000, 43000000, "Phạm Thanh Hoàng", $null, 0, 0, 0.02, 160000, 0, 0, 0, 0)
Set-TextValue $wsDebt 2 4 "08-05-2024"

Set-RowValues $wsDebt 3 1 @("Tổng", 1, 8000000, $null, $null, $null, $null, $null, $null, $null, 35000000, $null, 8000000, 43000000, 43000000, $null, $null, 0, 0, 0, 160000, 0, 0, 0, 0)

# ------------------------------------------------------------------
# Sheet 3: "Lương" (was sheet 2) - rebuilt from scratch with the new
# "Chiết khấu thu nợ" rows inserted + updated SÓC TRĂNG figures.
# ------------------------------------------------------------------
Set-RowValues $wsLuong 1 1 @("Danh mục lương", 13)
Set-RowValues $wsLuong 2 1 @("Tổng công tại CẦN THƠ", 0)
Set-RowValues $wsLuong 3 1 @("Lương công tác tại CẦN THƠ", 0)
Set-RowValues $wsLuong 4 1 @("Lương cơ bản tại CẦN THƠ", $null)
Set-RowValues $wsLuong 5 1 @("Chiết khấu sale chính tại CẦN THƠ", 0)
Set-RowValues $wsLuong 6 1 @("Chiết khấu sale phụ tại CẦN THƠ", 0)
Set-RowValues $wsLuong 7 1 @("Đơn 1 bác sĩ tại CẦN THƠ", 0)
Set-RowValues $wsLuong 8 1 @("Đơn 2 bác sĩ tại CẦN THƠ", 0)
Set-RowValues $wsLuong 9 1 @("Công phụ phẫu 1 tại CẦN THƠ", 0)
Set-RowValues $wsLuong 10 1 @("Công phụ phẫu 2 tại CẦN THƠ", 0)
Set-RowValues $wsLuong 11 1 @("Chiết khấu thu nợ tại CẦN THƠ", 0)
Set-RowValues $wsLuong 12 1 @("Ứng lương tại CẦN THƠ", 0)
Set-RowValues $wsLuong 13 1 @("Tổng công tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 14 1 @("Lương công tác tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 15 1 @("Lương cơ bản tại LONG XUYÊN", $null)
Set-RowValues $wsLuong 16 1 @("Chiết khấu sale chính tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 17 1 @("Chiết khấu sale phụ tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 18 1 @("Đơn 1 bác sĩ tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 19 1 @("Đơn 2 bác sĩ tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 20 1 @("Công phụ phẫu 1 tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 21 1 @("Công phụ phẫu 2 tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 22 1 @("Chiết khấu thu nợ tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 23 1 @("Ứng lương tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 24 1 @("Tổng công tại SÓC TRĂNG", 8)
Set-RowValues $wsLuong 25 1 @("Phụ cấp tại SÓC TRĂNG", 280000)
Set-RowValues $wsLuong 26 1 @("Lương cơ bản tại SÓC TRĂNG", 2571428.571428571)
Set-RowValues $wsLuong 27 1 @("Chiết khấu sale chính tại SÓC TRĂNG", 2560000)
Set-RowValues $wsLuong 28 1 @("Chiết khấu sale phụ tại SÓC TRĂNG", 0)
Set-RowValues $wsLuong 29 1 @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0)
Set-RowValues $wsLuong 30 1 @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0)
Set-RowValues $wsLuong 31 1 @("Công phụ phẫu 1 tại SÓC TRĂNG", 0)
Set-RowValues $wsLuong 32 1 @("Công phụ phẫu 2 tại SÓC TRĂNG", 0)
Set-RowValues $wsLuong 33 1 @("Chiết khấu thu nợ tại SÓC TRĂNG", 160000)
Set-RowValues $wsLuong 34 1 @("Ứng lương tại SÓC TRĂNG", 0)
Set-RowValues $wsLuong 35 1 @("Tổng lương tại CẦN THƠ", 0)
Set-RowValues $wsLuong 36 1 @("Tổng lương tại LONG XUYÊN", 0)
Set-RowValues $wsLuong 37 1 @("Tổng lương tại SÓC TRĂNG", 5571428.571428571)
Set-RowValues $wsLuong 38 1 @("Tổng lương tại HỆ THỐNG", 5571428.571428571)

$ws1.Select()
$ws1.Range("A1").Select()
